# Update the PRESUPUESTO (column G) values on the "VENTA MENSUAL" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G2").Value = 400
$ws.Range("G3").Value = 2500
$ws.Range("G4").Value = 200
$ws.Range("G5").Value = 2890
$ws.Range("G6").Value = 1000
$ws.Range("G11").Value = 500
$ws.Range("G12").Value = 7490
